# Update "想去人数" (F column) figures across sheets, matching the
# upstream scraper's refreshed output (gh-pages data update).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - 1st worksheet
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value  = 1704
$wsExhibit.Range("F3").Value  = 9195
$wsExhibit.Range("F4").Value  = 121
$wsExhibit.Range("F7").Value  = 1388
$wsExhibit.Range("F8").Value  = 210
$wsExhibit.Range("F9").Value  = 70
$wsExhibit.Range("F11").Value = 5977
$wsExhibit.Range("F15").Value = 4688
$wsExhibit.Range("F19").Value = 39
$wsExhibit.Range("F20").Value = 346
$wsExhibit.Range("F23").Value = 265
$wsExhibit.Range("F25").Value = 3109

# Sheet "演出" (Performances) - 2nd worksheet
$wsPerform = $wb.Worksheets.Item(2)
$wsPerform.Range("F2").Value = 50

# Sheet "全部类型" (All types) - 4th worksheet
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value  = 1704
$wsAll.Range("F3").Value  = 9195
$wsAll.Range("F4").Value  = 121
$wsAll.Range("F5").Value  = 50
$wsAll.Range("F8").Value  = 1388
$wsAll.Range("F9").Value  = 210
$wsAll.Range("F10").Value = 70
$wsAll.Range("F12").Value = 5977
$wsAll.Range("F16").Value = 4688
$wsAll.Range("F20").Value = 39
$wsAll.Range("F21").Value = 346
$wsAll.Range("F24").Value = 265
$wsAll.Range("F26").Value = 3109
